# Updated symbol list on Sun Dec 11 21:35:51 UTC 2022 with GitHub Actions
#
# This script applies the cell-level value changes described by the diff
# between the original "before.xlsx" snapshot of the cryptos worksheet and
# the refreshed scrape. The Price column (D) stores numeric-looking values
# as text, so those assignments are prefixed with a leading apostrophe to
# force a text literal instead of letting Excel auto-convert them to a
# floating point number (which would alter their printed representation,
# e.g. dropping significant trailing zeros).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 : BNB ---
$ws.Range("D2").Value = "'289.96"

# --- Row 4 : HuobiToken ---
$ws.Range("D4").Value = "'6.458"

# --- Row 5 : Cronos ---
$ws.Range("D5").Value = "'0.06390"

# --- Row 6 ---
$ws.Range("D6").Value = "'3.603"

# --- Row 8 ---
$ws.Range("D8").Value = "'6.588"

# --- Row 9 ---
$ws.Range("D9").Value = "'0.8280"

# --- Row 10 ---
$ws.Range("D10").Value = "'0.01427"

# --- Row 11 ---
$ws.Range("D11").Value = "'0.1689"

# --- Row 12 ---
$ws.Range("D12").Value = "'0.08790"

# --- Row 13 ---
$ws.Range("D13").Value = "'0.03678"

# --- Row 14 ---
$ws.Range("D14").Value = "'0.03203"

# --- Row 15 ---
$ws.Range("D15").Value = "'0.09197"

# --- Row 16 ---
$ws.Range("D16").Value = "'3.705"

# --- Row 17 ---
$ws.Range("D17").Value = "'0.001646"

# --- Row 18 ---
$ws.Range("D18").Value = "'0.04758"

# --- Row 19 ---
$ws.Range("D19").Value = "'0.006117"

# --- Row 20 ---
$ws.Range("D20").Value = "'0.006303"

# --- Row 21 ---
$ws.Range("D21").Value = "'0.001072"

# --- Row 23 ---
$ws.Range("D23").Value = "'3.779"

# --- Row 24 ---
$ws.Range("D24").Value = "'2.322"

# --- Row 26 ---
$ws.Range("D26").Value = "'0.1261"

# --- Row 27 : AAXToken -- Volume(1h) label lost its "Bestin24h" suffix ---
$ws.Range("E27").Value = "26AAXTokenAAB"

# --- Row 28 ---
$ws.Range("D28").Value = "'0.0002710"

# --- Row 40 ---
$ws.Range("D40").Value = "'0.04840"

# --- Row 41 : was BKEXToken, now KickToken (swapped with row 43) ---
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.007170"
$ws.Range("E41").Value = "40KickTokenKICK"

# --- Row 42 : CEJI ---
$ws.Range("D42").Value = "'0.004506"

# --- Row 43 : was KickToken, now BKEXToken (swapped with row 41) ---
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1118"
$ws.Range("E43").Value = "42BKEXTokenBKK"

# --- Row 44 : LocalTraders ---
$ws.Range("D44").Value = "'0.01185"

# --- Row 45 ---
$ws.Range("D45").Value = "'0.00006905"

# --- Row 46 ---
$ws.Range("D46").Value = "'0.00000000751"

# --- Row 47 ---
$ws.Range("D47").Value = "'0.9345"

# --- Row 48 : BOLO -- Volume(1h) label gained a "Bestin24h" suffix ---
$ws.Range("D48").Value = "'0.008636"
$ws.Range("E48").Value = "47BOLOBOLOBestin24h"

# --- Row 49 : CryptobidCoin ---
$ws.Range("D49").Value = "'0.00001903"

# --- Row 50 ---
$ws.Range("D50").Value = "'0.01242"
